$d = $word.ActiveDocument

# --- Hunk 1: split the run "{m" into two runs "{" and "m" -----------------
# Locate it robustly via Find so we do not depend on hard-coded offsets.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("{m", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found1) {
    # Collapse to just the "{" character (first char of the match).
    $brace = $d.Range($rng1.Start, $rng1.Start + 1)
    # Adding then immediately deleting a bookmark forces Word to break the
    # run at this boundary without touching any character formatting.
    $d.Bookmarks.Add("m2doc_tmp_split1", $brace)
    $d.Bookmarks("m2doc_tmp_split1").Delete()
}

# --- Hunk 2: split the run ".setWidth(100)}" into ".setWidth(100)" and "}" -
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(".setWidth(100)}", $true, $false, $false, `
                              $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # The "}" is the last character of the match.
    $closeBrace = $d.Range($rng2.End - 1, $rng2.End)

    # Step 1: break the run boundary (bookmark trick again) so the text
    # that stays behind (".setWidth(100)") keeps its original formatting
    # and identity instead of being silently re-merged with a neighbour.
    $d.Bookmarks.Add("m2doc_tmp_split2", $closeBrace)
    $d.Bookmarks("m2doc_tmp_split2").Delete()

    # Step 2: the newly isolated "}" run still carries the coloured
    # character formatting inherited from ".setWidth(100)". Replace its
    # content with a fresh plain run (no rPr) to match the reference
    # TokenIteratorFieldRewriterSplit output.
    $closeBrace2 = $d.Range($rng2.End - 1, $rng2.End)
    $closeBrace2.Delete()
    $insertion = $d.Range($rng2.End - 1, $rng2.End - 1)
    $insertion.InsertAfter("}")
}
